# The "Recorded By" column (G) lists the people/systems that recorded each
# attendance session as a comma-separated string, e.g.
#   "dnasr281@gmail.com, System"  or  "dnasr281@gmail.com, admin@admin.com"
# This sync normalises that list into alphabetical order, e.g.
#   "dnasr281@gmail.com, System"        -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System"       -> "System, backup@backdoor.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"
# Single-name cells (just "System" or just an email) and already-sorted /
# 3-name cells are left untouched, matching the upstream sync exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 157 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ([string]::IsNullOrEmpty($current)) { continue }

    $parts = $current -split ',\s*'

    if ($parts.Count -eq 2) {
        # Plain ordinal (case-sensitive) compare via .CompareTo() - "System"
        # (capital S) must sort before lowercase-leading emails, matching
        # the source data's byte-wise ordering. (The "-clt"/"-cle"
        # operators in this host are culture/case-INsensitive, so they
        # can't be used here.)
        if ($parts[0].CompareTo($parts[1]) -le 0) {
            $sorted = @($parts[0], $parts[1])
        } else {
            $sorted = @($parts[1], $parts[0])
        }
        $sortedJoined = [string]::Join(', ', $sorted)

        if ($sortedJoined -ne $current) {
            $cell.Value = $sortedJoined
        }
    }
}

Write-Host "Recorded-By column normalised to alphabetical order"
